# Insert two new weekly records for Brocoli / Femacal de La Calera just
# before the current row 365, shifting the rest of the table down by two
# rows (old A1:R417 -> new A1:R419).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("365:366").Insert()

# Row 365 - "Primera" quality record, dated 2021-10-22 (serial 44491)
$ws.Range("A365").Value = 3
$ws.Range("B365").Value = "Femacal de La Calera"
$ws.Range("C365").Value = "Coquimbo"
$ws.Range("D365").Value = 44491
$ws.Range("E365").Value = 5
$ws.Range("F365").Value = 100112023
$ws.Range("G365").Value = "Brócoli"
$ws.Range("H365").Value = "Sin especificar"
$ws.Range("I365").Value = "Primera"
$ws.Range("J365").Value = 2600
$ws.Range("K365").Value = 450
$ws.Range("L365").Value = 500
$ws.Range("M365").Value = 477
$ws.Range("N365").Value = "$/unidad"
$ws.Range("O365").Value = "Provincia de Quillota"
$ws.Range("P365").Value = 477
$ws.Range("Q365").Value = 1
$ws.Range("R365").Value = "Hortaliza"

# Row 366 - "Segunda" quality record, same date
$ws.Range("A366").Value = 3
$ws.Range("B366").Value = "Femacal de La Calera"
$ws.Range("C366").Value = "Coquimbo"
$ws.Range("D366").Value = 44491
$ws.Range("E366").Value = 5
$ws.Range("F366").Value = 100112023
$ws.Range("G366").Value = "Brócoli"
$ws.Range("H366").Value = "Sin especificar"
$ws.Range("I366").Value = "Segunda"
$ws.Range("J366").Value = 1600
$ws.Range("K366").Value = 600
$ws.Range("L366").Value = 600
$ws.Range("M366").Value = 600
$ws.Range("N366").Value = "$/unidad"
$ws.Range("O366").Value = "Provincia de Quillota"
$ws.Range("P366").Value = 600
$ws.Range("Q366").Value = 1
$ws.Range("R366").Value = "Hortaliza"
